$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column C (shifts existing C:Q to D:R)
$ws.Range("C1").EntireColumn.Insert()

# New header for the inserted column
$ws.Range("C1").Value = "Frame Openpose"

# New values for the inserted column
$ws.Range("C2").Formula = "=27*25"
$ws.Range("C3").Value = 1243
$ws.Range("C4").Value = 675
$ws.Range("C5").Value = 1243

# Update the active selection as recorded in the sheet view
$ws.Range("C4").Select()
